# Login.xlsx edit: add header/footer check sheets
# ------------------------------------------------------------------
# Adds two new worksheets - "LoginHeaderFooter" and
# "ForgotPasswordHeaderFooter" - each holding header/footer copy text used
# for a manual QA check, and leaves the new "ForgotPasswordHeaderFooter"
# sheet as the active tab (matching the reviewed workbook).

$wb = $excel.ActiveWorkbook

# ---- LoginHeaderFooter --------------------------------------------------
$sheet6 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$sheet6.Name = "LoginHeaderFooter"

# Write in the same order the original sheet's content was authored so the
# shared-string table gets built up with matching indices.
$sheet6.Range("A2").Value = "At vidit veritus duo, erat putent eu qui. Vim an numquam accumsan deserunt. Soluta delectus vim cu, ad nam agam epicuri democritum, ei torquatos scriptorem eum. Putent iracundia intellegat has an, mollis accusata scripserit pri ut. Exerci voluptua disputa"
$sheet6.Range("A1").Value = "Header Text"
$sheet6.Range("B1").Value = "Footer Text"
$sheet6.Range("B2").Value = "© 2014 - Reviewer Connect`nReviewer Connect Description"

$sheet6.Range("B2").WrapText = $true
$sheet6.Rows.Item(2).RowHeight = 135
$sheet6.Range("B2").Select()

# ---- ForgotPasswordHeaderFooter -----------------------------------------
$sheet7 = $wb.Worksheets.Add($null, $sheet6)
$sheet7.Name = "ForgotPasswordHeaderFooter"

$sheet7.Range("A2").Value = "Please enter your e-mail address in the space provided and click Send."
$sheet7.Range("A1").Value = "Header Text 1"
$sheet7.Range("B2").Value = "You will receive an e-mail with information for accessing your account."
$sheet7.Range("B1").Value = "Header Text 2"
$sheet7.Range("C1").Value = "Footer Text"
$sheet7.Range("C2").Value = "© 2014 - Reviewer Connect`nReviewer Connect Description"

$sheet7.Range("A2:C2").WrapText = $true
$sheet7.Rows.Item(2).RowHeight = 150

# Leave "ForgotPasswordHeaderFooter" as the active/selected sheet.
$sheet7.Activate()
$sheet7.Range("A1").Select()
